$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Attempting to (re-)populate the "Rating" column (D) for the first few
# books. Force Text formatting first so the values -- which carry stray
# leading/trailing newlines -- are kept verbatim instead of being
# auto-coerced into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "`n  3.91`n"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "`n  3.98`n"

$ws.Range("D4").Value = "Error"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "`n  4.27`n"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "`n  4.50`n"
